$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.510.47"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "2.488.44"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'314.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").Value = "'94.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.00%  "
$ws.Range("D7").Value = "'0.550"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("D10").Value = "'33.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "'0.110"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "2.869.80"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").Value = "'15.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("D16").Value = "2.469.26"
$ws.Range("E16").Value = "  -2.77%  "
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "41.463.12"
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").Value = "'6.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "'11.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.58%  "
$ws.Range("D22").Value = "'69.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").Value = "'236.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -3.46%  "
$ws.Range("D27").Value = "'24.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "'9.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").Value = "'37.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").Value = "'153.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("D33").Value = "'2.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0756"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "'17.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("D36").Value = "'3.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("E37").Value = "  -10.85%  "
$ws.Range("D38").Value = "'1.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("E40").Value = "  -4.57%  "
$ws.Range("D41").Value = "'4.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "'19.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.19%  "
$ws.Range("D44").Value = "1.986.94"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("E46").Value = "  -5.53%  "
$ws.Range("D47").Value = "'8.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").Value = "2.731.28"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'69.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("D50").Value = "'97.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").Value = "'0.178"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.26%  "
